$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank columns before old "contribuicoes" (M) for apoio_std/min/max
$ws.Columns("M:O").Insert()
# Insert 3 new blank columns before old "menor_ano" (now shifted to R) for contribuicoes_std/min/max
$ws.Columns("R:T").Insert()

# --- Header renames / new headers (row 1) ---
$ws.Range("H1").Value = "arrecadado_avg"
$ws.Range("I1").Value = "arrecadado_std"
$ws.Range("J1").Value = "arrecadado_min"
$ws.Range("K1").Value = "arrecadado_max"
$ws.Range("M1").Value = "apoio_std"
$ws.Range("N1").Value = "apoio_min"
$ws.Range("O1").Value = "apoio_max"
$ws.Range("Q1").Value = "contribuicoes_med"
$ws.Range("R1").Value = "contribuicoes_std"
$ws.Range("S1").Value = "contribuicoes_min"
$ws.Range("T1").Value = "contribuicoes_max"

# --- Recomputed data values per row (apoio_medio + new apoio_std/min/max, new contribuicoes_std/min/max) ---
# Row 2
$ws.Range("L2").Value = 76.06649705332148
$ws.Range("M2").Value = 30.35289755311455
$ws.Range("N2").Value = 30.69891629110986
$ws.Range("O2").Value = 216.9152091005833
$ws.Range("R2").Value = 410.5455139511234
$ws.Range("S2").Value = 14
$ws.Range("T2").Value = 3474

# Row 3
$ws.Range("L3").Value = 82.12933179093685
$ws.Range("M3").Value = 35.40469491555042
$ws.Range("N3").Value = 33.97203328414528
$ws.Range("O3").Value = 230.5829516876651
$ws.Range("R3").Value = 174.8393084425899
$ws.Range("S3").Value = 6
$ws.Range("T3").Value = 1540

# Row 4
$ws.Range("L4").Value = 84.02563974325884
$ws.Range("M4").Value = 39.83943189124467
$ws.Range("N4").Value = 11.52676430516467
$ws.Range("O4").Value = 254.2443749773306
$ws.Range("R4").Value = 155.4234048603088
$ws.Range("S4").Value = 1
$ws.Range("T4").Value = 1330

# Row 5
$ws.Range("L5").Value = 83.58617223970138
$ws.Range("M5").Value = 40.20112754081283
$ws.Range("N5").Value = 16.05922078302265
$ws.Range("O5").Value = 386.987447085286
$ws.Range("R5").Value = 248.1423584049538
$ws.Range("S5").Value = 2
$ws.Range("T5").Value = 2120

# Row 6
$ws.Range("L6").Value = 83.40563308192627
$ws.Range("M6").Value = 32.56399533953332
$ws.Range("N6").Value = 18.47818326605706
$ws.Range("O6").Value = 195.20880750356
$ws.Range("R6").Value = 501.9791222296157
$ws.Range("S6").Value = 1
$ws.Range("T6").Value = 7954

# Row 7
$ws.Range("L7").Value = 77.97569054482099
$ws.Range("M7").Value = 35.65744130420693
$ws.Range("N7").Value = 10.77163914429046
$ws.Range("O7").Value = 195.6882025465182
$ws.Range("R7").Value = 459.3640728169867
$ws.Range("S7").Value = 1
$ws.Range("T7").Value = 7954

# Row 8
$ws.Range("L8").Value = 77.73968420752422
$ws.Range("M8").Value = 41.30988410144955
$ws.Range("N8").Value = 12.19662302883409
$ws.Range("O8").Value = 247.2901437851162
$ws.Range("R8").Value = 685.9026443808731
$ws.Range("S8").Value = 2
$ws.Range("T8").Value = 7954

# Row 9
$ws.Range("L9").Value = 81.75853347173708
$ws.Range("M9").Value = 38.14520603080047
$ws.Range("N9").Value = 18.47818326605706
$ws.Range("O9").Value = 234.707661751482
$ws.Range("R9").Value = 508.76757749687
$ws.Range("S9").Value = 2
$ws.Range("T9").Value = 7954

# Row 10
$ws.Range("L10").Value = 82.92407682444032
$ws.Range("M10").Value = 38.77988334228132
$ws.Range("N10").Value = 11.93343625774652
$ws.Range("O10").Value = 230.5829516876651
$ws.Range("R10").Value = 247.246130032236
$ws.Range("S10").Value = 1
$ws.Range("T10").Value = 2684

# Row 11
$ws.Range("L11").Value = 71.78666858221021
$ws.Range("M11").Value = 29.97870848948209
$ws.Range("N11").Value = 16.18065842403185
$ws.Range("O11").Value = 216.9152091005833
$ws.Range("R11").Value = 331.1306730481258
$ws.Range("S11").Value = 3
$ws.Range("T11").Value = 3474

# Row 12
$ws.Range("L12").Value = 71.78666858221021
$ws.Range("M12").Value = 29.97870848948209
$ws.Range("N12").Value = 16.18065842403185
$ws.Range("O12").Value = 216.9152091005833
$ws.Range("R12").Value = 331.1306730481258
$ws.Range("S12").Value = 3
$ws.Range("T12").Value = 3474

# Row 13
$ws.Range("L13").Value = 79.12786981308152
$ws.Range("M13").Value = 35.35639160943987
$ws.Range("N13").Value = 20.51363271354002
$ws.Range("O13").Value = 233.3973531230909
$ws.Range("R13").Value = 328.5928536530323
$ws.Range("S13").Value = 1
$ws.Range("T13").Value = 4584

# Row 14
$ws.Range("L14").Value = 84.32898346466456
$ws.Range("M14").Value = 39.18079463334893
$ws.Range("N14").Value = 11.52676430516467
$ws.Range("O14").Value = 195.6882025465182
$ws.Range("R14").Value = 147.3429131999576
$ws.Range("S14").Value = 1
$ws.Range("T14").Value = 612

# Row 15
$ws.Range("L15").Value = 83.48264574282582
$ws.Range("M15").Value = 40.20970245451376
$ws.Range("N15").Value = 21.00493274015408
$ws.Range("O15").Value = 247.2901437851162
$ws.Range("R15").Value = 174.0035582476068
$ws.Range("S15").Value = 2
$ws.Range("T15").Value = 770

# Row 16
$ws.Range("L16").Value = 83.64941179158359
$ws.Range("M16").Value = 36.99022540033587
$ws.Range("N16").Value = 16.18065842403185
$ws.Range("O16").Value = 254.2443749773306
$ws.Range("R16").Value = 183.1544319258093
$ws.Range("S16").Value = 1
$ws.Range("T16").Value = 1540

# Row 17
$ws.Range("L17").Value = 70.93306185876429
$ws.Range("M17").Value = 30.19803349932243
$ws.Range("N17").Value = 20.33774597757668
$ws.Range("O17").Value = 159.7763429092917
$ws.Range("R17").Value = 975.3935739169402
$ws.Range("S17").Value = 3
$ws.Range("T17").Value = 7954

# Row 18
$ws.Range("L18").Value = 76.84102373029619
$ws.Range("M18").Value = 33.37177139781743
$ws.Range("N18").Value = 16.18065842403185
$ws.Range("O18").Value = 226.5579622472015
$ws.Range("R18").Value = 457.1568742729124
$ws.Range("S18").Value = 1
$ws.Range("T18").Value = 7954

# Row 19
$ws.Range("L19").Value = 79.19230719197579
$ws.Range("M19").Value = 27.01736191709247
$ws.Range("N19").Value = 40.63189862969614
$ws.Range("O19").Value = 130.9739254174068
$ws.Range("R19").Value = 143.8139060050855
$ws.Range("S19").Value = 2
$ws.Range("T19").Value = 467

# Row 20
$ws.Range("L20").Value = 78.67160937524555
$ws.Range("M20").Value = 39.22105185666557
$ws.Range("N20").Value = 14.90596347946683
$ws.Range("O20").Value = 461.5197709071476
$ws.Range("R20").Value = 402.3874992420548
$ws.Range("S20").Value = 1
$ws.Range("T20").Value = 7954

# Row 21
$ws.Range("L21").Value = 69.51944033042635
$ws.Range("M21").Value = 28.22546328675653
$ws.Range("N21").Value = 17.82064921105857
$ws.Range("O21").Value = 196.4212117364618
$ws.Range("R21").Value = 448.9315388299796
$ws.Range("S21").Value = 4
$ws.Range("T21").Value = 4584

# Row 22
$ws.Range("L22").Value = 74.35971583315494
$ws.Range("M22").Value = 27.7214279498591
$ws.Range("N22").Value = 17.83984513748501
$ws.Range("O22").Value = 156.0426904908593
$ws.Range("R22").Value = 190.7586913141665
$ws.Range("S22").Value = 1
$ws.Range("T22").Value = 1879
